$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

$ws.Range("B2").Value = 17150
$ws.Range("B3").Value = 16650
$ws.Range("B4").Value = 15950
$ws.Range("B5").Value = 15950
